$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74: Gonçalves
$ws.Range("C74").Value = 4355

# Row 77: Guaxupé
$ws.Range("C77").Value = 52078
$ws.Range("D77").Value = 49430
$ws.Range("E77").Value = 172.59
$ws.Range("F77").Value = 3128709
$ws.Range("G77").Value = "guaxupeano"
$ws.Range("H77").Value = "HEBER HAMILTON QUINTELLA"
$ws.Range("I77").Value = 3

# Row 78: Heliodora
$ws.Range("C78").Value = 6591
$ws.Range("D78").Value = 6121
$ws.Range("E78").Value = 39.76
$ws.Range("F78").Value = 3129202
$ws.Range("G78").Value = "heliodorense"
$ws.Range("H78").Value = "ALEX LEOPOLDINO DE LIMA"
$ws.Range("I78").Value = 1

# Row 79: Ibiraci
$ws.Range("C79").Value = 13986
$ws.Range("D79").Value = 12176
$ws.Range("E79").Value = 21.66
$ws.Range("F79").Value = 3129707
$ws.Range("G79").Value = "ibiraciense"
$ws.Range("H79").Value = "ISMAEL SILVA CÂNDIDO"
$ws.Range("I79").Value = 1

# Row 80: Ibitiúra de Minas
$ws.Range("C80").Value = 3492
$ws.Range("D80").Value = 3382
$ws.Range("E80").Value = 49.51
$ws.Range("F80").Value = 3129905
$ws.Range("G80").Value = "ibitiurense"
$ws.Range("H80").Value = "ALEXANDRE DE CASSIO BORGES"
$ws.Range("I80").Value = 1

# Row 81: Ibituruna
$ws.Range("C81").Value = 2996
$ws.Range("D81").Value = 2866
$ws.Range("E81").Value = 18.72
$ws.Range("F81").Value = 3130002
$ws.Range("G81").Value = "ibiturunense"
$ws.Range("H81").Value = "FRANCISCO ANTONIO PEREIRA"
$ws.Range("I81").Value = 1

# Row 82: Ijaci
$ws.Range("C82").Value = 6610
$ws.Range("D82").Value = 5859
$ws.Range("E82").Value = 55.67
$ws.Range("F82").Value = 3130408
$ws.Range("G82").Value = "ijaciense"
$ws.Range("H82").Value = "FABIANO DA SILVA MORETI"
$ws.Range("I82").Value = 1

# Row 83: Ilicínea
$ws.Range("C83").Value = 12444
$ws.Range("D83").Value = 11488
$ws.Range("E83").Value = 30.53
$ws.Range("F83").Value = 3130507
$ws.Range("G83").Value = "ilicineaense"
$ws.Range("H83").Value = "NIRLEI CRISTIANI"
$ws.Range("I83").Value = 1

# Row 84: Inconfidentes
$ws.Range("C84").Value = 7358
$ws.Range("D84").Value = 6908
$ws.Range("E84").Value = 46.17
$ws.Range("F84").Value = 3130606
$ws.Range("G84").Value = "inconfidentino"
$ws.Range("H84").Value = "ROSÂNGELA MARIA DANTAS"

# Row 85: Ingaí
$ws.Range("C85").Value = 2776
$ws.Range("D85").Value = 2629
$ws.Range("E85").Value = 8.6
$ws.Range("F85").Value = 3130804
$ws.Range("G85").Value = "ingaiense"
$ws.Range("H85").Value = "GIULLIANO RIBEIRO PINTO"

# Row 86: Ipuiúna
$ws.Range("C86").Value = 10118
$ws.Range("D86").Value = 9521
$ws.Range("E86").Value = 31.93
$ws.Range("F86").Value = 3131505
$ws.Range("G86").Value = "ipuiunense"
$ws.Range("H86").Value = "ELDER CASSIO DE SOUZA OLIVA"

# Row 87: Itajubá
$ws.Range("C87").Value = 97334
$ws.Range("D87").Value = 90658
$ws.Range("E87").Value = 307.49
$ws.Range("F87").Value = 3132404
$ws.Range("G87").Value = "Itajubense"
$ws.Range("H87").Value = "CHRISTIAN GONÇALVES TIBURZIO E SILVA"
$ws.Range("I87").Value = 3

# Row 88: Itamogi
$ws.Range("C88").Value = 10157
$ws.Range("D88").Value = 10349
$ws.Range("E88").Value = 42.47
$ws.Range("F88").Value = 3132909
$ws.Range("G88").Value = "itamogiense"
$ws.Range("H88").Value = "RONALDO PEREIRA DIAS"

# Row 89: Itamonte
$ws.Range("C89").Value = 15714
$ws.Range("D89").Value = 14003
$ws.Range("E89").Value = 32.43
$ws.Range("F89").Value = 3133006
$ws.Range("G89").Value = "itamontense"
$ws.Range("H89").Value = "ALEXANDRE AUGUSTO MOREIRA SANTOS"
$ws.Range("I89").Value = 2

# Row 90: Itanhandu
$ws.Range("C90").Value = 15423
$ws.Range("D90").Value = 14175
$ws.Range("E90").Value = 98.87
$ws.Range("F90").Value = 3133105
$ws.Range("G90").Value = "itanhanduense"
$ws.Range("H90").Value = "CARLOS GONÇALVES DA FONSECA"
$ws.Range("I90").Value = 2

# Row 91: Itapeva
$ws.Range("C91").Value = 9881
$ws.Range("D91").Value = 8664
$ws.Range("E91").Value = 48.85
$ws.Range("F91").Value = 3133600
$ws.Range("G91").Value = "itapevense"
$ws.Range("H91").Value = "DANIEL PEREIRA DO COUTO"

# Row 92: Itaú de Minas
$ws.Range("C92").Value = 16199
$ws.Range("D92").Value = 14945
$ws.Range("E92").Value = 97.41
$ws.Range("F92").Value = 3133758
$ws.Range("G92").Value = "itauense"
$ws.Range("H92").Value = "NORIVAL FRANCISCO DE LIMA"
$ws.Range("I92").Value = 2

# Row 93: Itumirim
$ws.Range("C93").Value = 6000
$ws.Range("D93").Value = 6139
$ws.Range("E93").Value = 26.15
$ws.Range("F93").Value = 3134301
$ws.Range("G93").Value = "itumirinense"
$ws.Range("H93").Value = "CARLOS ALBERTO NASCIMENTO"

# Row 94: Jacutinga
$ws.Range("C94").Value = 26264
$ws.Range("D94").Value = 22772
$ws.Range("E94").Value = 65.48
$ws.Range("F94").Value = 3134905
$ws.Range("G94").Value = "jacutinguense"
$ws.Range("H94").Value = "MELQUIADES DE ARAUJO"
$ws.Range("I94").Value = 2

# Row 95: Jacuí
$ws.Range("C95").Value = 7691
$ws.Range("D95").Value = 7502
$ws.Range("E95").Value = 18.33
$ws.Range("F95").Value = 3134806
$ws.Range("G95").Value = "jacuiense"
$ws.Range("H95").Value = "MARIA CONCEICAO DOS REIS PEREIRA"

# Row 96: Jesuânia
$ws.Range("C96").Value = 4780
$ws.Range("D96").Value = 4768
$ws.Range("E96").Value = 30.99
$ws.Range("F96").Value = 3135902
$ws.Range("G96").Value = "jesuanense"
$ws.Range("H96").Value = "JOSÉ LAÉRCIO BRANDÃO DE CASTRO"

# Row 97: Juruaia
$ws.Range("C97").Value = 10681
$ws.Range("D97").Value = 9238
$ws.Range("E97").Value = 41.92
$ws.Range("F97").Value = 3136900
$ws.Range("G97").Value = "juruaiense"
$ws.Range("H97").Value = "ÁLVARO MARIANO JÚNIOR"

# Row 98: Lambari
$ws.Range("C98").Value = 20907
$ws.Range("D98").Value = 19554
$ws.Range("E98").Value = 91.76
$ws.Range("F98").Value = 3137809
$ws.Range("G98").Value = "lambariense"
$ws.Range("H98").Value = "MARCELO GIOVANI DE SOUSA"
$ws.Range("I98").Value = 2

# Row 99: Lavras
$ws.Range("C99").Value = 104783
$ws.Range("D99").Value = 92200
$ws.Range("E99").Value = 163.26
$ws.Range("F99").Value = 3138203
$ws.Range("G99").Value = "lavrense"
$ws.Range("H99").Value = "JUSSARA MENICUCCI DE OLIVEIRA"
$ws.Range("I99").Value = 3

# Row 100: Liberdade
$ws.Range("C100").Value = 5031
$ws.Range("D100").Value = 5346
$ws.Range("E100").Value = 13.32
$ws.Range("F100").Value = 3138500
$ws.Range("G100").Value = "libertense"
$ws.Range("H100").Value = "WALTER DE ASSIS TOLEDO JUNIOR"

# Row 101: Luminárias
$ws.Range("C101").Value = 5438
$ws.Range("D101").Value = 5422
$ws.Range("E101").Value = 10.84
$ws.Range("F101").Value = 3138708
$ws.Range("G101").Value = "luminarense"
$ws.Range("H101").Value = "ECIO CARVALHO REZENDE"
